# Natmi following Dr Hou advice
# Rewrite the LR-pair data rows: row 2 gets new values, and three more
# sending-cluster rows (FAPs, M2, sCs -> now ECs/FAPs/M2/sCs) are added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 1) are unchanged.

# Row 2: ECs -> Tnfsf13b -> Tnfrsf13c -> sCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13c"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 1.429365
$ws.Range("H2").Value = 4.288095
$ws.Range("I2").Value = 0.327844155149115
$ws.Range("J2").Value = 0.327844155149115
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.832688
$ws.Range("N2").Value = 2.498064
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.19021508312
$ws.Range("R2").Value = 10.71193574808
$ws.Range("S2").Value = 0.327844155149115
$ws.Range("T2").Value = 0.327844155149115

# Row 3: FAPs -> Tnfsf13b -> Tnfrsf13c -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13c"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.643401666666666
$ws.Range("H3").Value = 4.930204999999999
$ws.Range("I3").Value = 0.3769363535408946
$ws.Range("J3").Value = 0.3769363535408946
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.832688
$ws.Range("N3").Value = 2.498064
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.368440847013333
$ws.Range("R3").Value = 12.31596762312
$ws.Range("S3").Value = 0.3769363535408946
$ws.Range("T3").Value = 0.3769363535408946

# Row 4: M2 -> Tnfsf13b -> Tnfrsf13c -> sCs
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13c"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.256465666666666
$ws.Range("H4").Value = 3.769397
$ws.Range("I4").Value = 0.2881873593953978
$ws.Range("J4").Value = 0.2881873593953979
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.832688
$ws.Range("N4").Value = 2.498064
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.046243883045333
$ws.Range("R4").Value = 9.416194947407998
$ws.Range("S4").Value = 0.2881873593953978
$ws.Range("T4").Value = 0.2881873593953979

# Row 5: sCs -> Tnfsf13b -> Tnfrsf13c -> sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Tnfsf13b"
$ws.Range("C5").Value = "Tnfrsf13c"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03065933333333333
$ws.Range("H5").Value = 0.091978
$ws.Range("I5").Value = 0.00703213191459268
$ws.Range("J5").Value = 0.007032131914592681
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.832688
$ws.Range("N5").Value = 2.498064
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.02552965895466667
$ws.Range("R5").Value = 0.229766930592
$ws.Range("S5").Value = 0.00703213191459268
$ws.Range("T5").Value = 0.007032131914592681
